# New crime data collected - weekly CompStat update (113th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -78.947368421052

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = 43.478260869565
$ws.Range("N15").Value = -54.794520547945

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -55.555555555555
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = -18.666666666666
$ws.Range("L16").Value = -20.261437908496
$ws.Range("M16").Value = -61.993769470405
$ws.Range("N16").Value = -88.544600938967

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -24.242424242424
$ws.Range("I17").Value = 360
$ws.Range("J17").Value = 372
$ws.Range("K17").Value = -3.225806451612
$ws.Range("L17").Value = -6.735751295336
$ws.Range("M17").Value = 9.422492401215
$ws.Range("N17").Value = -49.790794979079

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -11.111111111111
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = -1.666666666666
$ws.Range("L18").Value = -20.805369127516
$ws.Range("M18").Value = -69.350649350649
$ws.Range("N18").Value = -88.793922127255

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -70
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -35.135135135135
$ws.Range("I19").Value = 407
$ws.Range("J19").Value = 331
$ws.Range("K19").Value = 22.960725075528
$ws.Range("L19").Value = 9.703504043126
$ws.Range("M19").Value = -33.821138211382
$ws.Range("N19").Value = -89.082618025751

# --- Row 20 (G.L.A.) ---
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -22.727272727272
$ws.Range("I20").Value = 217
$ws.Range("J20").Value = 172
$ws.Range("K20").Value = 26.162790697674
$ws.Range("L20").Value = 13.020833333333
$ws.Range("M20").Value = -13.545816733067
$ws.Range("N20").Value = -87.197640117994

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -4
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = -30.4
$ws.Range("I21").Value = 1265
$ws.Range("J21").Value = 1177
$ws.Range("K21").Value = 7.476635514018
$ws.Range("L21").Value = -1.937984496124
$ws.Range("M21").Value = -34.860968074150
$ws.Range("N21").Value = -84.884693511769

# --- Row 23 (Housing) ---
# F23 switches from a numeric cell to the shared "0" text placeholder;
# copy style+text from G23 (which already holds that placeholder) then done.
$ws.Range("G23").Copy($ws.Range("F23"))
$ws.Range("L23").Value = -50
$ws.Range("M23").Value = -57.142857142857

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 112.5
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 5.617977528089
$ws.Range("I24").Value = 1219
$ws.Range("J24").Value = 899
$ws.Range("K24").Value = 35.59510567297
$ws.Range("L24").Value = 27.911857292759
$ws.Range("M24").Value = 20.335636722606

# --- Row 25 (Misd. Assault) ---
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 58.823529411764
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 454
$ws.Range("K25").Value = 10.132158590308
$ws.Range("L25").Value = 7.526881720430
$ws.Range("M25").Value = -33.065595716198

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = 1
# G26: numeric 1 -> shared "0" text placeholder (copy style+text from C26)
$ws.Range("C26").Copy($ws.Range("G26"))
# H26: numeric 100 -> shared "***.*" text placeholder (copy style+text from E26)
$ws.Range("E26").Copy($ws.Range("H26"))
$ws.Range("L26").Value = 0

# --- Row 27 (Other Sex Crimes) ---
# C27/D27 switch from the "0" placeholder to real numbers; E27 switches from
# "***.*" to a real percentage. Borrow number styles from row 16 (same column
# style family: s=15 for counts, s=16 for percents) then overwrite the value.
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D16").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("E16").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = 12.765957446808
$ws.Range("L27").Value = 130.434782608696

# --- Row 28 (Shooting Vic.) ---
# D28: numeric 2 -> shared "0" placeholder; E28: numeric -100 -> shared "***.*"
$ws.Range("C28").Copy($ws.Range("D28"))
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("M28").Value = -45.614035087719
$ws.Range("N28").Value = -79.605263157894

# --- Row 29 (Shooting Inc.) ---
$ws.Range("C29").Copy($ws.Range("D29"))
$ws.Range("E22").Copy($ws.Range("E29"))
$ws.Range("M29").Value = -46.808510638297
$ws.Range("N29").Value = -81.617647058823
